$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2490.625
$ws.Range("I43").Value = 1770
$ws.Range("K43").Value = 1770
$ws.Range("M43").Value = -1701
$ws.Range("H64").Value = 8000
$ws.Range("I64").Value = 8000
$ws.Range("K64").Value = 8000
$ws.Range("M64").Value = -7752
$ws.Range("H67").Value = 8000
$ws.Range("I67").Value = 8000
$ws.Range("K67").Value = 8000
$ws.Range("M67").Value = -7142
$ws.Range("H74").Value = 7282.5
$ws.Range("I74").Value = 7282.5
$ws.Range("K74").Value = 7282.5
$ws.Range("M74").Value = -6346.5
$ws.Range("H77").Value = 7282.5
$ws.Range("I77").Value = 7282.5
$ws.Range("K77").Value = 36412.5
$ws.Range("M77").Value = -31732.5
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 34343.5
$ws.Range("I63").Value = 53018.332
$ws.Range("K63").Value = 53018.332
$ws.Range("M63").Value = -52332.332
$ws.Range("H66").Value = 34343.5
$ws.Range("I66").Value = 53018.332
$ws.Range("K66").Value = 265091.66
$ws.Range("M66").Value = -261659.66
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H122").Value = 3799.6
$ws.Range("I122").Value = 3799.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11398.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8948.799999999999
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H124").Value = 100000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 100000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820
$ws.Range("H125").Value = 100715
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 100715
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 100715
$ws.Range("N125").Value = -110555
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 55000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 55000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 55000
$ws.Range("N127").Value = -64920
$ws.Range("H128").Value = 99999.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 99999.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 99999.5
$ws.Range("N128").Value = -109959.5
$ws.Range("H129").Value = 44753.332
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 44753.332
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 44753.332
$ws.Range("N129").Value = -54753.332
$ws.Range("H130").Value = 70099.25
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 70099.25
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 70099.25
$ws.Range("N130").Value = -80139.25
$ws.Range("H131").Value = 95107.836
$ws.Range("I131").Value = 50650
$ws.Range("J131").Value = 103999.4
$ws.Range("K131").Value = 50650
$ws.Range("L131").Value = 103999.4
$ws.Range("M131").Value = -45610
$ws.Range("N131").Value = -114079.4
$ws.Range("H132").Value = 3941
$ws.Range("I132").Value = 4102.4
$ws.Range("J132").Value = 3779.6
$ws.Range("K132").Value = 12307.2
$ws.Range("L132").Value = 11338.8
$ws.Range("M132").Value = -9777.199999999999
$ws.Range("N132").Value = -16398.8
$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 98333.336
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 98333.336
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 98333.336
$ws.Range("N135").Value = -108473.336
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 150000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 150000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280
$ws.Range("H139").Value = 100000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 100000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
$ws.Range("H140").Value = 80000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 80000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 383.33334
$ws.Range("H62").Value = 25335.334
$ws.Range("J62").Value = 25335.334
$ws.Range("L62").Value = 25335.334
$ws.Range("N62").Value = -26583.334
$ws.Range("H65").Value = 25335.334
$ws.Range("J65").Value = 25335.334
$ws.Range("L65").Value = 126676.67
$ws.Range("N65").Value = -132916.67
$ws.Range("H124").Value = 105163
$ws.Range("J124").Value = 105163
$ws.Range("L124").Value = 105163
$ws.Range("N124").Value = -110073
$ws.Range("H134").Value = 8240.4
$ws.Range("I134").Value = 9750.5
$ws.Range("K134").Value = 29251.5
$ws.Range("M134").Value = -26716.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 125
$ws.Range("I7").Value = 125
$ws.Range("K7").Value = 375
$ws.Range("M7").Value = -263
$ws.Range("H56").Value = 12333.333
$ws.Range("I56").Value = 12333.333
$ws.Range("K56").Value = 12333.333
$ws.Range("M56").Value = -11803.333
$ws.Range("H68").Value = 4001.5
$ws.Range("J68").Value = 4001.5
$ws.Range("L68").Value = 12004.5
$ws.Range("N68").Value = -13626.5
$ws.Range("H71").Value = 4001.5
$ws.Range("J71").Value = 4001.5
$ws.Range("L71").Value = 36013.5
$ws.Range("N71").Value = -44125.5
$ws.Range("H81").Value = 6399.8
$ws.Range("J81").Value = 6250
$ws.Range("L81").Value = 18750
$ws.Range("N81").Value = -20996
$ws.Range("H84").Value = 6399.8
$ws.Range("J84").Value = 6250
$ws.Range("L84").Value = 56250
$ws.Range("N84").Value = -67482
$ws.Range("H92").Value = 447.16666
$ws.Range("I92").Value = 395
$ws.Range("K92").Value = 1185
$ws.Range("M92").Value = 63

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5999
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 5999
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H80").Value = 10312.5
$ws.Range("I80").Value = 4400
$ws.Range("K80").Value = 4400
$ws.Range("M80").Value = -3402
$ws.Range("H83").Value = 10312.5
$ws.Range("I83").Value = 4400
$ws.Range("K83").Value = 22000
$ws.Range("M83").Value = -17008
$ws.Range("H132").Value = 4055.7144
$ws.Range("I132").Value = 3297.3333
$ws.Range("J132").Value = 4624.5
$ws.Range("K132").Value = 9891.999899999999
$ws.Range("L132").Value = 13873.5
$ws.Range("M132").Value = -7361.999899999999
$ws.Range("N132").Value = -18933.5
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12591.923
$ws.Range("I22").Value = 15782.667
$ws.Range("J22").Value = 9857
$ws.Range("K22").Value = 15782.667
$ws.Range("L22").Value = 9857
$ws.Range("M22").Value = -15487.667
$ws.Range("N22").Value = -10447
$ws.Range("H27").Value = 12591.923
$ws.Range("I27").Value = 15782.667
$ws.Range("J27").Value = 9857
$ws.Range("K27").Value = 15782.667
$ws.Range("L27").Value = 9857
$ws.Range("M27").Value = -15675.667
$ws.Range("N27").Value = -10071
$ws.Range("H40").Value = 8214.143
$ws.Range("I40").Value = 8099.8
$ws.Range("J40").Value = 8500
$ws.Range("K40").Value = 8099.8
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = -7963.8
$ws.Range("N40").Value = -8772
$ws.Range("H46").Value = 5333.3335
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -5812
$ws.Range("N46").Value = -4376

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 89998.336
$ws.Range("J130").Value = 89998.336
$ws.Range("L130").Value = 89998.336
$ws.Range("N130").Value = -100038.336
